# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
# The "Periodo Mora" column (E16:E20) holds the period code "2507" for every
# worker row on the statement; this update rolls it forward to "2508" for the
# new account statement (part 1).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16:E20").Value = "2508"
